{"js": "// The numbered \"steps\" list is rewritten:\n//   - paragraphs 1-3 (\"1. npm init ...\", \"2. npm i express\",\n//     \"3. touch server.js ...\") get their text retyped as a single run\n//     each (this clears the spell-check proofErr markers and the\n//     run-splitting that npm/package.json/etc. produced), and\n//   - a new 4th step (\"4. \u0417\u0430\u043f\u0443\u0441\u043a - node server\") is appended at the end.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0 is the youtube link and stays untouched; paragraphs 1..3 are\n// the three existing numbered steps that need their runs collapsed into a\n// single clean run (same visible text, no proofErr spell-check markers).\nconst newStepText = [\n  \"1. npm init - create a package.json\",\n  \"2. npm i express\",\n  \"3. touch server.js - \u0441\u043e\u0437\u0434\u0430\u043d\u0438\u0435 \u0444\u0430\u0439\u043b\u0430 server.js\",\n];\n\nfor (let i = 0; i < newStepText.length; i++) {\n  const oldParagraph = paragraphs.items[i + 1];\n  // Inserting a fresh paragraph right before the old one copies its\n  // paragraph/run formatting (the en-US language mark) but starts with a\n  // clean single run and no proofErr markers; then drop the old paragraph.\n  oldParagraph.insertParagraph(newStepText[i], \"Before\");\n  oldParagraph.delete();\n  await context.sync();\n}\n\n// Re-load paragraphs and add the new 4th step after the last one (old\n// paragraph 3, \"3. touch server.js - ...\").\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = refreshed.items[refreshed.items.length - 1];\nlastParagraph.insertParagraph(\"4. \u0417\u0430\u043f\u0443\u0441\u043a - node server\", \"After\");\nawait context.sync();\n", "ps1": "# The numbered \"steps\" list is rewritten:\n#   - paragraphs 2-4 (\"1. npm init ...\", \"2. npm i express\",\n#     \"3. touch server.js ...\") get their text retyped as a single run\n#     each (this clears the spell-check proofErr markers and the\n#     run-splitting that npm/package.json/etc. produced), and\n#   - a new 4th step (\"4. \u0417\u0430\u043f\u0443\u0441\u043a - node server\") is appended at the end.\n$d = $word.ActiveDocument\n\n$newStepText = @(\n    \"1. npm init - create a package.json\",\n    \"2. npm i express\",\n    \"3. touch server.js - \u0441\u043e\u0437\u0434\u0430\u043d\u0438\u0435 \u0444\u0430\u0439\u043b\u0430 server.js\"\n)\n\n# Paragraphs(1) is the youtube link and stays untouched; Paragraphs(2..4)\n# are the three existing numbered steps that need their runs collapsed\n# into a single clean run (same visible text, no proofErr markers).\nfor ($i = 0; $i -lt $newStepText.Length; $i++) {\n    $idx = $i + 2\n    $oldRange = $d.Paragraphs($idx).Range\n    # Inserting a fresh empty paragraph right before the old one copies\n    # its paragraph formatting (the en-US language mark); filling that\n    # empty paragraph's Range.Text gives one clean run with no proofErr\n    # markers. Then the old (still messy) paragraph is deleted.\n    $oldRange.InsertParagraphBefore()\n    $d.Paragraphs($idx).Range.Text = $newStepText[$i]\n    $d.Paragraphs($idx + 1).Range.Delete()\n}\n\n# Append the new 4th step after the last paragraph (the rewritten\n# \"3. touch server.js - ...\" paragraph).\n$lastIndex = $d.Paragraphs.Count\n$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"4. \u0417\u0430\u043f\u0443\u0441\u043a - node server\"\n"}
